$d = $word.ActiveDocument

function ReplaceOnce($old, $new) {
    $rng = $d.Content
    $null = $rng.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

# 1. "elicitate" -> "elicit" (counter measure sentence)
ReplaceOnce "easier to elicitate the counter measure using the tactics." "easier to elicit the counter measure using the tactics."

# 2. Merge the two trailing runs of the "Some cases may not..." sentence (no wording change)
ReplaceOnce "Some cases may not have the clear distinction of source, target, and boundary conditions and so applying the operators may not provide the best resolution." "Some cases may not have the clear distinction of source, target, and boundary conditions and so applying the operators may not provide the best resolution."

# 3. Merge the runs of "False or incorrect data can severely affect the output" (no wording change)
ReplaceOnce "False or incorrect data can severely affect the output" "False or incorrect data can severely affect the output"

# 4. Merge the runs of "If the source, target and boundary conditions ..." (no wording change)
ReplaceOnce "If the source, target and boundary conditions are identified, this technique can be efficient to come with a good resolution." "If the source, target and boundary conditions are identified, this technique can be efficient to come with a good resolution."

# 5. "huge an random" -> "huge a random"
ReplaceOnce "Good point to start the search of risk.  If the feature set it huge an random this technique can provide a good starting point." "Good point to start the search of risk.  If the feature set it huge a random this technique can provide a good starting point."

# 6. "If the there are" -> "If there are"
ReplaceOnce "If the there are more well-defined components in the project, component inspection would be a better choice than this technique" "If there are more well-defined components in the project, component inspection would be a better choice than this technique"

# 7. "coflicts" -> "conflicts"
ReplaceOnce "Identifying Defects and coflicts" "Identifying Defects and conflicts"

# 8. "Conflics" -> "Conflicts"
ReplaceOnce "Documenting Conflics" "Documenting Conflicts"

# 9. "Identifcation" -> "Identification"
ReplaceOnce "Risk Identifcation" "Risk Identification"

# 10. Split "Appendix A to Appendix E" at the point Word left its _GoBack bookmark
$rng = $d.Content
$null = $rng.Find.Execute(", please see Appendix A to Appendix E.")
$splitPoint = $rng.Start + 17
$bmRange = $d.Range($splitPoint, $splitPoint)
$d.Bookmarks.Add("_GoBack", $bmRange)
